$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.054.55"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "1.643.35"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("D5").Value = "'215.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").Value = "'0.5052"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("D7").Value = "'1.010"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "'0.2580"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").Value = "'0.06437"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").Value = "'0.07725"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "1.645.41"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").Value = "'4.256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "1.870.75"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "'0.5459"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "0.0₅7922"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "'63.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "26.049.02"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "'1.010"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'204.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").Value = "'4.303"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").Value = "'10.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").Value = "'5.970"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").Value = "'1.011"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "'1.947"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.29%  "
$ws.Range("D26").Value = "'141.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("D27").Value = "'0.1157"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").Value = "'15.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "'6.755"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.03%  "
$ws.Range("D30").Value = "'0.05069"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").Value = "'3.252"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.57%  "
$ws.Range("D33").Value = "'3.198"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("E34").Value = "  -1.24%  "
$ws.Range("D35").Value = "'2.343"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").Value = "'0.8971"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.87%  "
$ws.Range("D37").Value = "'2.621"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.76%  "
$ws.Range("D38").Value = "'0.5639"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").Value = "1.146.76"
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("D41").Value = "'2.567"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'1.010"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").Value = "'5.672"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "'0.8155"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("D45").Value = "'99.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "1.781.01"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").Value = "0.0₈112"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "'0.4535"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("D49").Value = "'1.011"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "'54.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").Value = "'0.05042"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.98%  "
